# Update cryptos list (Price/Volume(1h) columns) to the latest scraped values.
# Column D (Price) values are forced to text with a leading apostrophe so that
# Excel does not reinterpret number-looking strings (e.g. "60.947.02",
# "7.80") as numeric values and strip formatting / change cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''60.947.02'
$ws.Range("E2").Value = '  -2.04%  '
$ws.Range("D3").Value = '''2.413.38'
$ws.Range("E3").Value = '  -1.41%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '''567.52'
$ws.Range("E5").Value = '  -2.82%  '
$ws.Range("D6").Value = '''138.44'
$ws.Range("E6").Value = '  -3.30%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").Value = '''0.526'
$ws.Range("E8").Value = '  -1.21%  '
$ws.Range("D9").Value = '''2.396.64'
$ws.Range("E9").Value = '  -1.91%  '
$ws.Range("E10").Value = '  -2.76%  '
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("E12").Value = '  -3.05%  '
$ws.Range("D13").Value = '''0.337'
$ws.Range("E13").Value = '  -2.09%  '
$ws.Range("D14").Value = '''25.89'
$ws.Range("E14").Value = '  -2.41%  '
$ws.Range("D15").Value = '''2.856.71'
$ws.Range("D16").Value = '''0.0000170'
$ws.Range("E16").Value = '  -3.53%  '
$ws.Range("D17").Value = '''60.763.34'
$ws.Range("E17").Value = '  -2.09%  '
$ws.Range("D18").Value = '''2.404.95'
$ws.Range("E18").Value = '  -1.15%  '
$ws.Range("E19").Value = '  +7.76%  '
$ws.Range("E20").Value = '  -1.93%  '
$ws.Range("D21").Value = '''321.94'
$ws.Range("E23").Value = '  +1.82%  '
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("E25").Value = '  -5.78%  '
$ws.Range("D26").Value = '''64.78'
$ws.Range("E26").Value = '  -1.51%  '
$ws.Range("D27").Value = '''576.79'
$ws.Range("E27").Value = '  -4.43%  '
$ws.Range("D28").Value = '''8.15'
$ws.Range("E28").Value = '  -10.56%  '
$ws.Range("D30").Value = '''0.0₃0915'
$ws.Range("E30").Value = '  -5.57%  '
$ws.Range("D31").Value = '''7.80'
$ws.Range("E31").Value = '  -2.50%  '
$ws.Range("E32").Value = '  -6.14%  '
$ws.Range("E33").Value = '  -4.66%  '
$ws.Range("E34").Value = '  -3.10%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").Value = '''152.17'
$ws.Range("E36").Value = '  -1.11%  '
$ws.Range("E37").Value = '  -2.38%  '
$ws.Range("E38").Value = '  -7.01%  '
$ws.Range("E39").Value = '  -3.04%  '
$ws.Range("D40").Value = '''18.17'
$ws.Range("E40").Value = '  -1.49%  '
$ws.Range("D41").Value = '''5.09'
$ws.Range("E41").Value = '  -3.81%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '''1.66'
$ws.Range("E43").Value = '  -3.20%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '''41.08'
$ws.Range("E44").Value = '  -4.90%  '
$ws.Range("E45").Value = '  -10.07%  '
$ws.Range("D46").Value = '''141.92'
$ws.Range("E46").Value = '  -0.40%  '
$ws.Range("D47").Value = '''0.0₆0264'
$ws.Range("D48").Value = '''3.49'
$ws.Range("E48").Value = '  -4.13%  '
$ws.Range("E49").Value = '  -2.87%  '
$ws.Range("D50").Value = '''0.0499'
$ws.Range("E50").Value = '  -4.23%  '
$ws.Range("D51").Value = '''19.21'
$ws.Range("E51").Value = '  -3.24%  '
